# Refresh of the "쿼리1" query-table sheet: new query result data (ranking
# reshuffled, values updated, refresh timestamp bumped) plus the leftover
# "applyNumberFormat"-only style getting cleared off column B so those
# cells fall back to the default style, and the selection left wherever
# the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("쿼리1")

# New query results (rank, BJ name, monthly cumulative balloons) plus the
# bumped "last refreshed" timestamp for every data row.
$refreshed = 46025.637322974537

$data = @(
    @{ Row = 2;  Name = "태영";   Value = 127144 },
    @{ Row = 3;  Name = "으냉이"; Value = 80731 },
    @{ Row = 4;  Name = "하루묭"; Value = 57673 },
    @{ Row = 5;  Name = "우리밍"; Value = 52173 },
    @{ Row = 6;  Name = "빵지니"; Value = 48927 },
    @{ Row = 7;  Name = "한쪼니"; Value = 46331 },
    @{ Row = 8;  Name = "임밍지"; Value = 33510 },
    @{ Row = 9;  Name = "윤하랑"; Value = 17690 },
    @{ Row = 10; Name = "히요코"; Value = 2154 },
    @{ Row = 11; Name = "하랑e"; Value = $null }
)

foreach ($item in $data) {
    $r = $item.Row

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $item.Name
    $bCell.ClearFormats()

    $cCell = $ws.Cells.Item($r, 3)
    if ($null -eq $item.Value) {
        $cCell.ClearContents()
    } else {
        $cCell.Value = $item.Value
    }

    $ws.Cells.Item($r, 4).Value = $refreshed
}

# Leftover selection from the editing session.
[void]$ws.Range("H12").Select()
